$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'328.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.25%"
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'44.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'0.73%"
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'5.481"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-1.11%"
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'0.08042"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'0.37%"
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'2.023"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'6.76%"
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.9530"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'1.05%"
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.1115"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-5.99%"
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.1874"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'2.05%"
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'10.16"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'1.15%"
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.09952"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'3.20%"
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.04741"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'6.61%"
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'0.1063"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.30%"
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'0.001263"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-2.15%"
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.04098"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-2.04%"
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'0.005863"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-1.24%"
$ws.Range("E16").Style = "Normal"

$ws.Range("D18").Value = "'4.412"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'3.50%"
$ws.Range("E18").Style = "Normal"

$ws.Range("E19").Value = "'2.42%"
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'0.3407"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-1.06%"
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'0.1401"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-1.15%"
$ws.Range("E21").Style = "Normal"

$ws.Range("E22").Value = "'2.80%"
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'0.001311"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'5.05%"
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'0.004352"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'1.56%"
$ws.Range("E24").Style = "Normal"

$ws.Range("E25").Value = "'-0.75%"
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'0.0003745"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-6.35%"
$ws.Range("E26").Style = "Normal"

$ws.Range("D38").Value = "'0.02572"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'-2.19%"
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'0.05657"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'3.07%"
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.007738"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'2.10%"
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.1397"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'0.48%"
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.007357"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-10.09%"
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'0.002010"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'0.42%"
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.008496"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-3.45%"
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.00007084"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-0.18%"
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'0.06%"
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.0005808"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.05%"
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'0.003504"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'6.90%"
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'0.003503"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'53.92%"
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.00002103"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.06%"
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'0.0002003"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.06%"
$ws.Range("E51").Style = "Normal"

